# Ready for Nov 6 class
# ---------------------------------------------------------------------------
# 1) Slide 2 ("Previously, in IMM120" / "Midterms were finished!") gets new
#    bullet content about Collisions, including two hyperlinked runs.
# 2) The presentation theme's colour scheme ("Blue II") is swapped for the
#    "Median" palette (same slot layout, new RGB values).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- Slide 2: add the new "Collisions" bullets -----------------------------
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$tr.Text = "Midterms were finished!`r" + `
           "`r" + `
           "Collisions`r" + `
           "Guide Posted`r" + `
           "Just realized how poorly the guides look and work, sorry!`r" + `
           "https://github.com/crhallberg/IMM120 `r"

# Paragraph 4 - "Guide Posted" (level 2 / lvl="1") with a hyperlink
$pGuide = $tr.Paragraphs(4)
$pGuide.IndentLevel = 2
$pGuide.ActionSettings(1).Hyperlink.Address = "https://crhallberg.github.io/IMM120/"

# Paragraph 5 - apology line (level 2 / lvl="1")
$pSorry = $tr.Paragraphs(5)
$pSorry.IndentLevel = 2

# Paragraph 6 - repo link (level 2 / lvl="1"); only the URL text itself is
# hyperlinked, the trailing space stays plain.
$pLink = $tr.Paragraphs(6)
$pLink.IndentLevel = 2
$urlChars = $pLink.Characters(1, 36)
$urlChars.ActionSettings(1).Hyperlink.Address = "https://github.com/crhallberg/IMM120"

# --- Theme: swap "Blue II" colours for the "Median" palette -----------------
$scheme = $s.ThemeColorScheme
$scheme.Item(3).RGB  = 5594999    # dk2      -> 775F55
$scheme.Item(4).RGB  = 12836331   # lt2      -> EBDDC3
$scheme.Item(5).RGB  = 13809300   # accent1  -> 94B6D2
$scheme.Item(6).RGB  = 4686045    # accent2  -> DD8047
$scheme.Item(7).RGB  = 8498085    # accent3  -> A5AB81
$scheme.Item(8).RGB  = 6075096    # accent4  -> D8B25C
$scheme.Item(9).RGB  = 10332027   # accent5  -> 7BA79D
$scheme.Item(10).RGB = 9211030    # accent6  -> 968C8C
$scheme.Item(11).RGB = 1423095    # hlink    -> F7B615
$scheme.Item(12).RGB = 279664     # folHlink -> 704404
